$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$templates = $wb.Worksheets.Item("templates")

# Populate the new header row on Sheet1.
# Column H ("Email Date") is written before column G ("Email Sent") so that
# the shared-string table ends up with the same insertion order as the
# original authoring session.
$ws1.Range("A1").Value = "First Name"
$ws1.Range("B1").Value = "Last Name"
$ws1.Range("C1").Value = "Location"
$ws1.Range("D1").Value = "Email"
$ws1.Range("E1").Value = "Child's Name"
$ws1.Range("F1").Value = "Child Age"
$ws1.Range("H1").Value = "Email Date"
$ws1.Range("G1").Value = "Email Sent"
$ws1.Range("I1").Value = "Response"
$ws1.Range("J1").Value = "Response Date"

# Approximate the "best fit" column widths from the finished sheet.
$ws1.Columns.Item(1).ColumnWidth = 9.498697916666666
$ws1.Columns.Item(2).ColumnWidth = 9.166666666666666
$ws1.Columns.Item(3).ColumnWidth = 7.166666666666667
$ws1.Columns.Item(4).ColumnWidth = 4.998697916666667
$ws1.Columns.Item(5).ColumnWidth = 11.166666666666666
$ws1.Columns.Item(6).ColumnWidth = 7.998697916666667
$ws1.Columns.Item(7).ColumnWidth = 9.166666666666666
$ws1.Columns.Item(8).ColumnWidth = 9.330729166666666
$ws1.Columns.Item(9).ColumnWidth = 8.166666666666666
$ws1.Columns.Item(10).ColumnWidth = 12.498697916666666

# Update the remembered selection on the templates sheet before it loses focus.
[void]$templates.Range("G3").Select()

# Move Sheet1 in front of templates so it becomes the first/active tab.
$ws1.Move($templates)

# Re-fetch the worksheet reference by name: after Move() the old object
# reference tracks the worksheet that now occupies the original slot,
# not the worksheet that actually moved.
$ws1 = $wb.Worksheets.Item("Sheet1")

# Restore the active selection on the (now active) Sheet1 tab.
[void]$ws1.Range("M9").Select()
